$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1281.3125
$ws.Range("I19").Value = 1044.7142
$ws.Range("J19").Value = 1465.3334
$ws.Range("K19").Value = 1044.7142
$ws.Range("L19").Value = 1465.3334
$ws.Range("M19").Value = -869.7141999999999
$ws.Range("N19").Value = -1815.3334
$ws.Range("H98").Value = 2729.4583
$ws.Range("I98").Value = 2377.5908
$ws.Range("J98").Value = 6600
$ws.Range("K98").Value = 2377.5908
$ws.Range("L98").Value = 6600
$ws.Range("M98").Value = -879.5907999999999
$ws.Range("N98").Value = -9596
$ws.Range("H122").Value = 2729.4583
$ws.Range("I122").Value = 2377.5908
$ws.Range("J122").Value = 6600
$ws.Range("K122").Value = 7132.7724
$ws.Range("L122").Value = 19800
$ws.Range("M122").Value = -4682.7724
$ws.Range("N122").Value = -24700
$ws.Range("H137").Value = 26673.36
$ws.Range("I137").Value = 945.6087
$ws.Range("J137").Value = 63657
$ws.Range("K137").Value = 2836.8261
$ws.Range("L137").Value = 190971
$ws.Range("M137").Value = -286.8261000000002
$ws.Range("N137").Value = -196071
$ws.Range("H138").Value = 1837.1208
$ws.Range("I138").Value = 1537.0317
$ws.Range("J138").Value = 2512.3215
$ws.Range("K138").Value = 4611.0951
$ws.Range("L138").Value = 7536.9645
$ws.Range("M138").Value = 528.9048999999995
$ws.Range("N138").Value = -17816.9645
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 347861.12
$ws.Range("I2").Value = 427965.62
$ws.Range("K2").Value = 427965.62
$ws.Range("M2").Value = -427852.62
$ws.Range("H32").Value = 2159.99
$ws.Range("I32").Value = 1876.3656
$ws.Range("J32").Value = 5928.143
$ws.Range("K32").Value = 1876.3656
$ws.Range("L32").Value = 5928.143
$ws.Range("M32").Value = -1589.3656
$ws.Range("N32").Value = -6502.143
$ws.Range("H45").Value = 1461.6842
$ws.Range("I45").Value = 1100.1428
$ws.Range("J45").Value = 1672.5834
$ws.Range("K45").Value = 1100.1428
$ws.Range("L45").Value = 1672.5834
$ws.Range("M45").Value = -723.1428000000001
$ws.Range("N45").Value = -2426.5834
$ws.Range("H110").Value = 183.54546
$ws.Range("I110").Value = 121.375
$ws.Range("K110").Value = 121.375
$ws.Range("M110").Value = 1923.625
$ws.Range("H116").Value = 347861.12
$ws.Range("I116").Value = 427965.62
$ws.Range("K116").Value = 427965.62
$ws.Range("M116").Value = -425671.62
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 347861.12
$ws.Range("I3").Value = 427965.62
$ws.Range("K3").Value = 427965.62
$ws.Range("M3").Value = -427851.62
$ws.Range("H80").Value = 6828.0625
$ws.Range("I80").Value = 382.66666
$ws.Range("J80").Value = 10695.3
$ws.Range("K80").Value = 382.66666
$ws.Range("L80").Value = 10695.3
$ws.Range("M80").Value = 615.33334
$ws.Range("N80").Value = -12691.3
$ws.Range("H83").Value = 6828.0625
$ws.Range("I83").Value = 382.66666
$ws.Range("J83").Value = 10695.3
$ws.Range("K83").Value = 1913.3333
$ws.Range("L83").Value = 53476.5
$ws.Range("M83").Value = 3078.6667
$ws.Range("N83").Value = -63460.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2267.8
$ws.Range("I31").Value = 1918.4166
$ws.Range("J31").Value = 3665.3333
$ws.Range("K31").Value = 1918.4166
$ws.Range("L31").Value = 3665.3333
$ws.Range("M31").Value = -1623.4166
$ws.Range("N31").Value = -4255.3333
$ws.Range("H34").Value = 2267.8
$ws.Range("I34").Value = 1918.4166
$ws.Range("J34").Value = 3665.3333
$ws.Range("K34").Value = 1918.4166
$ws.Range("L34").Value = 3665.3333
$ws.Range("M34").Value = -1716.4166
$ws.Range("N34").Value = -4069.3333
$ws.Range("H58").Value = 926193.4
$ws.Range("I58").Value = 1403379.1
$ws.Range("J58").Value = 1646
$ws.Range("K58").Value = 1403379.1
$ws.Range("L58").Value = 1646
$ws.Range("M58").Value = -1403176.1
$ws.Range("N58").Value = -2052
$ws.Range("H99").Value = 1880.5
$ws.Range("I99").Value = 1539.1
$ws.Range("K99").Value = 1539.1
$ws.Range("M99").Value = -41.09999999999991
$ws.Range("H126").Value = 1880.5
$ws.Range("I126").Value = 1539.1
$ws.Range("K126").Value = 4617.299999999999
$ws.Range("M126").Value = -2147.299999999999
$ws.Range("H134").Value = 1556.0834
$ws.Range("I134").Value = 1417.7709
$ws.Range("K134").Value = 4253.3127
$ws.Range("M134").Value = -1718.3127
$ws.Range("H136").Value = 926193.4
$ws.Range("I136").Value = 1403379.1
$ws.Range("J136").Value = 1646
$ws.Range("K136").Value = 4210137.300000001
$ws.Range("L136").Value = 4938
$ws.Range("M136").Value = -4207587.300000001
$ws.Range("N136").Value = -10038
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 160
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H131").Value = 14346.394
$ws.Range("J131").Value = 16807.768
$ws.Range("L131").Value = 50423.304
$ws.Range("N131").Value = -60503.304
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3009.8
$ws.Range("I70").Value = 2667.7
$ws.Range("K70").Value = 2667.7
$ws.Range("M70").Value = -2397.7
$ws.Range("H73").Value = 3009.8
$ws.Range("I73").Value = 2667.7
$ws.Range("K73").Value = 2667.7
$ws.Range("M73").Value = -1731.7
$ws.Range("H102").Value = 3194.6875
$ws.Range("I102").Value = 3660.1428
$ws.Range("J102").Value = 2306.0908
$ws.Range("K102").Value = 3660.1428
$ws.Range("L102").Value = 2306.0908
$ws.Range("M102").Value = -2038.1428
$ws.Range("N102").Value = -5550.0908
$ws.Range("H132").Value = 535693.9399999999
$ws.Range("I132").Value = 785991.0600000001
$ws.Range("J132").Value = 2452.1738
$ws.Range("K132").Value = 2357973.18
$ws.Range("L132").Value = 7356.5214
$ws.Range("M132").Value = -2355443.18
$ws.Range("N132").Value = -12416.5214
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 16773.834
$ws.Range("I16").Value = 16773.834
$ws.Range("K16").Value = 16773.834
$ws.Range("M16").Value = -16603.834
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1151.9828
$ws.Range("I132").Value = 637.1818
$ws.Range("J132").Value = 2769.9285
$ws.Range("K132").Value = 1911.5454
$ws.Range("L132").Value = 8309.7855
$ws.Range("M132").Value = 618.4546
$ws.Range("N132").Value = -13369.7855
$ws.Range("H136").Value = 12921737
$ws.Range("I136").Value = 14621711
$ws.Range("J136").Value = 1938.4
$ws.Range("K136").Value = 43865133
$ws.Range("L136").Value = 5815.200000000001
$ws.Range("M136").Value = -43862583
$ws.Range("N136").Value = -10915.2
